# EEPROM setup over serial PoC
# Update the "stock" sheet: remove four now-unused stock rows (C4/47u,
# C8/470u, J2 barrel jack, L1 inductor), shrinking the parts list, and add a
# new stock line for Q3 (BC547 transistor), while fixing up a handful of
# "# needed" (K column) quantities elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stock")
$ws.Activate() | Out-Null

# --- 1) Quantity corrections (K column) -----------------------------------
$ws.Cells.Item(5, 11).Value2 = 5
$ws.Cells.Item(8, 11).Value2 = 5
$ws.Cells.Item(18, 11).Value2 = 4
$ws.Cells.Item(19, 11).Value2 = 5

# --- 2) Remove obsolete stock rows -----------------------------------------
# Delete bottom-up so the row numbers of the still-to-delete rows don't shift.
$ws.Rows(46).Delete() | Out-Null   # L1   Inductor_THT ... 470u
$ws.Rows(45).Delete() | Out-Null   # J2   Barrel_Jack_Switch
$ws.Rows(42).Delete() | Out-Null   # C8   470u electrolytic
$ws.Rows(41).Delete() | Out-Null   # C4   47u electrolytic

# Rows 43..53 have now shifted up to 41..49, in order:
#   41 D3 D5 / BZX84C3V3
#   42 D9    / CUS10S30
#   43 Q2 Q8 / BSS84
#   44 Q4 Q5 Q7 Q9 Q10 / BC847B
#   45 Q6    / DMP4065
#   46 R2 R4 R14 R40 / 22k
#   47 R36   / 6k8
#   48 R6 R10 R11 / 4k7
#   49 (#N/A formula) / 0R

# Grab the existing "Arwill" supplier shared string (untouched row 5) before
# we start writing the new row; re-typing the accented text risks mangling
# the shared-string bytes differently than the original mojibake. Also grab
# the Hyperlink cell style already used by the other manual-link (J) cells
# so the new one reuses the same style record instead of minting a new one.
$arwill = $ws.Cells.Item(5, 9).Value2
$hyperlinkStyle = $ws.Cells.Item(5, 10).Style

# --- 3) Append the new Q3 (BC547) stock row at row 50 ----------------------
# Fill cells in the same order the shared-string table expects them to have
# been entered (Reference/Value/Footprint, then MPN/Manufacturer, then the
# product-page hyperlink, and finally the datasheet link).
$newRow = 50
$ws.Cells.Item($newRow, 1).Value2 = "Q3 "
$ws.Cells.Item($newRow, 2).Value2 = 1
$ws.Cells.Item($newRow, 3).Value2 = "BC547"
$ws.Cells.Item($newRow, 4).Value2 = "Package_TO_SOT_THT:TO-92_Inline_Wide"
$ws.Cells.Item($newRow, 7).Value2 = "BC547B"
$ws.Cells.Item($newRow, 8).Value2 = "CDIL"
$ws.Cells.Item($newRow, 9).Value2 = $arwill
$ws.Cells.Item($newRow, 11).Value2 = 8

$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 10), "https://arwill.hu/termekek/felvezetok/tranzisztorok/bc/bc547b-tranzisztor-160303/") | Out-Null
$ws.Cells.Item($newRow, 10).Style = $hyperlinkStyle

$ws.Cells.Item($newRow, 5).Value2 = "https://arwill.hu/forras/termek/felvezetok/tranzisztorok/bc/bc547b-tranzisztor-160303.pdf"

# --- 4) Fix up the autofilter range + named range to the new extents -------
$newRange = $ws.Range("A1:K50")
$newRange.AutoFilter() | Out-Null
$newRange.AutoFilter() | Out-Null

for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "stock!_FilterDatabase") {
        $n.RefersTo = "=stock!`$A`$1:`$K`$50"
    }
}

# --- 5) Move the selection, matching the author's last cursor position -----
$ws.Range("K51").Select() | Out-Null
